$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.823114002679745
$ws.Range("C2").Value = 0.4813958968679231
$ws.Range("D2").Value = 0.03834743138354213
$ws.Range("E2").Value = 0.4165535847911457
$ws.Range("F2").Value = 1.607424994274538
$ws.Range("N2").Value = 1.212118694768108
$ws.Range("B3").Value = 1.630456767043654
$ws.Range("C3").Value = 0.4204967157134547
$ws.Range("D3").Value = 0.03753683670579733
$ws.Range("E3").Value = 0.3628436459713242
$ws.Range("F3").Value = 1.525055899423933
$ws.Range("N3").Value = 1.220976500462896
$ws.Range("B4").Value = 1.5130977584638
$ws.Range("C4").Value = 0.3832437444423817
$ws.Range("D4").Value = 0.0370738913139661
$ws.Range("E4").Value = 0.3300359615345201
$ws.Range("F4").Value = 1.475737788692058
$ws.Range("N4").Value = 1.226962336557378
$ws.Range("B5").Value = 1.465500177147646
$ws.Range("C5").Value = 0.368095284847584
$ws.Range("D5").Value = 0.03689378231912599
$ws.Range("E5").Value = 0.3167055049957952
$ws.Range("F5").Value = 1.455950422279628
$ws.Range("N5").Value = 1.229537901306273
$ws.Range("B6").Value = 1.457610138403879
$ws.Range("C6").Value = 0.3655817777130892
$ws.Range("D6").Value = 0.03686438576837503
$ws.Range("E6").Value = 0.3144942316036321
$ws.Range("F6").Value = 1.452683308102579
$ws.Range("N6").Value = 1.229973765968055
$ws.Range("B7").Value = 1.512454930879699
$ws.Range("C7").Value = 0.383039319144018
$ws.Range("D7").Value = 0.03707142796038454
$ws.Range("E7").Value = 0.3298560297946125
$ws.Range("F7").Value = 1.475469681206945
$ws.Range("N7").Value = 1.226996521383782
$ws.Range("B8").Value = 1.756487451233681
$ws.Range("C8").Value = 0.4603673259097718
$ws.Range("D8").Value = 0.03806060385504395
$ws.Range("E8").Value = 0.397996633269031
$ws.Range("F8").Value = 1.57875967946643
$ws.Range("N8").Value = 1.215058535366232
$ws.Range("B9").Value = 2.24280884307052
$ws.Range("C9").Value = 0.6132401623679584
$ws.Range("D9").Value = 0.04028529937923508
$ws.Range("E9").Value = 0.5331574849092533
$ws.Range("F9").Value = 1.791563452685381
$ws.Range("N9").Value = 1.196047647162132
$ws.Range("B10").Value = 2.605402809847362
$ws.Range("C10").Value = 0.7264963117532375
$ws.Range("D10").Value = 0.04210642825474054
$ws.Range("E10").Value = 0.6336659378606555
$ws.Range("F10").Value = 1.954569570134367
$ws.Range("N10").Value = 1.184843868109169
$ws.Range("B11").Value = 2.771630732378753
$ws.Range("C11").Value = 0.7782661748723285
$ws.Range("D11").Value = 0.04297831604792179
$ws.Range("E11").Value = 0.6797117557207599
$ws.Range("F11").Value = 2.030262665696227
$ws.Range("N11").Value = 1.180365045838158
$ws.Range("B12").Value = 2.834771016760556
$ws.Range("C12").Value = 0.7979092261151663
$ws.Range("D12").Value = 0.04331496359854015
$ws.Range("E12").Value = 0.697199565656021
$ws.Range("F12").Value = 2.059154535319635
$ws.Range("N12").Value = 1.178759399911769
$ws.Range("B13").Value = 2.821163880238942
$ws.Range("C13").Value = 0.7936769618805215
$ws.Range("D13").Value = 0.04324216841251882
$ws.Range("E13").Value = 0.6934308976683923
$ws.Range("F13").Value = 2.052921877308734
$ws.Range("N13").Value = 1.179101159600208
$ws.Range("B14").Value = 2.776821393528394
$ws.Range("C14").Value = 0.779881422245353
$ws.Range("D14").Value = 0.04300588095929214
$ws.Range("E14").Value = 0.6811494339380175
$ws.Range("F14").Value = 2.03263499648321
$ws.Range("N14").Value = 1.180231127950719
$ws.Range("B15").Value = 2.749685795068274
$ws.Range("C15").Value = 0.7714364191049867
$ws.Range("D15").Value = 0.04286199935777546
$ws.Range("E15").Value = 0.6736334985529879
$ws.Range("F15").Value = 2.020238668069311
$ws.Range("N15").Value = 1.180935083537165
$ws.Range("B16").Value = 2.5945659884934
$ws.Range("C16").Value = 0.7231183056020996
$ws.Range("D16").Value = 0.04205034256276718
$ws.Range("E16").Value = 0.6306636515447366
$ws.Range("F16").Value = 1.949654445585793
$ws.Range("N16").Value = 1.185149133291389
$ws.Range("B17").Value = 2.499739728419115
$ws.Range("C17").Value = 0.6935426905603208
$ws.Range("D17").Value = 0.04156371071912446
$ws.Range("E17").Value = 0.6043893134583556
$ws.Range("F17").Value = 1.90675285851168
$ws.Range("N17").Value = 1.18789364674906
$ws.Range("B18").Value = 2.445318296630035
$ws.Range("C18").Value = 0.6765548321291135
$ws.Range("D18").Value = 0.04128788803190275
$ws.Range("E18").Value = 0.5893072213077204
$ws.Range("F18").Value = 1.882221601045472
$ws.Range("N18").Value = 1.18953027760891
$ws.Range("B19").Value = 2.426912495739657
$ws.Range("C19").Value = 0.6708069353937844
$ws.Range("D19").Value = 0.04119519203637623
$ws.Range("E19").Value = 0.5842057328244863
$ws.Range("F19").Value = 1.873940364122973
$ws.Range("N19").Value = 1.190094337879785
$ws.Range("B20").Value = 2.509821648676223
$ws.Range("C20").Value = 0.6966886363742901
$ws.Range("D20").Value = 0.04161509007005293
$ws.Range("E20").Value = 0.6071831013337743
$ws.Range("F20").Value = 1.911304790056136
$ws.Range("N20").Value = 1.187595468608365
$ws.Range("B21").Value = 2.789840535537792
$ws.Range("C21").Value = 0.7839324249746369
$ws.Range("D21").Value = 0.04307510647793578
$ws.Range("E21").Value = 0.6847553686611576
$ws.Range("F21").Value = 2.038587486369266
$ws.Range("N21").Value = 1.179896763329893
$ws.Range("B22").Value = 2.97398057198302
$ws.Range("C22").Value = 0.841179685161535
$ws.Range("D22").Value = 0.04406720988179558
$ws.Range("E22").Value = 0.7357543503820239
$ws.Range("F22").Value = 2.123109523238497
$ws.Range("N22").Value = 1.175392796021939
$ws.Range("B23").Value = 2.875594996329028
$ws.Range("C23").Value = 0.8106037976108951
$ws.Range("D23").Value = 0.04353415888736123
$ws.Range("E23").Value = 0.7085061186113251
$ws.Range("F23").Value = 2.077873984331234
$ws.Range("N23").Value = 1.177747856566313
$ws.Range("B24").Value = 2.505263315905779
$ws.Range("C24").Value = 0.6952663060501436
$ws.Range("D24").Value = 0.04159184918428593
$ws.Range("E24").Value = 0.6059199573793705
$ws.Range("F24").Value = 1.90924644670568
$ws.Range("N24").Value = 1.187730091995817
$ws.Range("B25").Value = 2.110357201227885
$ws.Range("C25").Value = 0.5717327513650616
$ws.Range("D25").Value = 0.03965163600347665
$ws.Range("E25").Value = 0.4964009163727496
$ws.Range("F25").Value = 1.732855016866665
$ws.Range("N25").Value = 1.200711276772694
